$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 397 ("「創造的であり続けるための40の方法」...") entirely.
# This shifts every row below it up by one (398 -> 397, ..., 589 -> 588),
# matching the diff which shows the dimension shrinking from C589 to C588.
$ws.Rows("397:397").Delete()
